# Insert one new weekly data row for "Membrillo" (Vega Modelo de Temuco) above the
# current row 280. This shifts the existing rows 280..336 down to 281..337 and
# leaves a single blank row (with the date column's style already inherited)
# at row 280, which we then populate with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 280..336 down to 281..337, creating a new blank row 280.
$ws.Rows.Item(280).Insert()

$r = 280

$ws.Cells.Item($r, 1).Value = 10
$ws.Cells.Item($r, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item($r, 3).Value = "La Araucanía"
$ws.Cells.Item($r, 4).Value = 45180
$ws.Cells.Item($r, 5).Value = 9
$ws.Cells.Item($r, 6).Value = "Fruta"
$ws.Cells.Item($r, 7).Value = 100104
$ws.Cells.Item($r, 8).Value = "Frutos de pepita"
$ws.Cells.Item($r, 9).Value = 100104003
$ws.Cells.Item($r, 10).Value = "Membrillo"
$ws.Cells.Item($r, 11).Value = "Champion"
$ws.Cells.Item($r, 12).Value = "Primera"
$ws.Cells.Item($r, 13).Value = 90
$ws.Cells.Item($r, 14).Value = 15000
$ws.Cells.Item($r, 15).Value = 15000
$ws.Cells.Item($r, 16).Value = 15000
$ws.Cells.Item($r, 17).Value = "$/bandeja 18 kilos granel"
$ws.Cells.Item($r, 18).Value = "Región de O'Higgins"
$ws.Cells.Item($r, 19).Value = 833
$ws.Cells.Item($r, 20).Value = 18
